$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1547.3077
$ws.Range("I28").Value = 554.375
$ws.Range("J28").Value = 3136
$ws.Range("K28").Value = 554.375
$ws.Range("L28").Value = 3136
$ws.Range("M28").Value = -69.375
$ws.Range("N28").Value = -4106
$ws.Range("H132").Value = 99797.28999999999
$ws.Range("I132").Value = 114763.5
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 344290.5
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -341760.5
$ws.Range("N132").Value = -35060
$ws.Range("H137").Value = 3000
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3703.3333
$ws.Range("I32").Value = 3069.375
$ws.Range("J32").Value = 8775
$ws.Range("K32").Value = 3069.375
$ws.Range("L32").Value = 8775
$ws.Range("M32").Value = -2782.375
$ws.Range("N32").Value = -9349
$ws.Range("H61").Value = 2190.6667
$ws.Range("I61").Value = 2190.6667
$ws.Range("K61").Value = 2190.6667
$ws.Range("M61").Value = -1978.6667
$ws.Range("H76").Value = 29879.5
$ws.Range("J76").Value = 29879.5
$ws.Range("L76").Value = 29879.5
$ws.Range("N76").Value = -30555.5
$ws.Range("H79").Value = 29879.5
$ws.Range("J79").Value = 29879.5
$ws.Range("L79").Value = 29879.5
$ws.Range("N79").Value = -32219.5
$ws.Range("H124").Value = 81819
$ws.Range("J124").Value = 81819
$ws.Range("L124").Value = 81819
$ws.Range("N124").Value = -91639
$ws.Range("H132").Value = 2842.4
$ws.Range("I132").Value = 2842.4
$ws.Range("K132").Value = 8527.200000000001
$ws.Range("M132").Value = -5997.200000000001
$ws.Range("H135").Value = 15000
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = ""
$ws.Range("H136").Value = 2190.6667
$ws.Range("I136").Value = 2190.6667
$ws.Range("K136").Value = 6572.000100000001
$ws.Range("M136").Value = -4022.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 34037.418
$ws.Range("I107").Value = 40596.1
$ws.Range("K107").Value = 40596.1
$ws.Range("M107").Value = -38676.1
$ws.Range("H134").Value = 500.8
$ws.Range("I134").Value = 500.8
$ws.Range("K134").Value = 1502.4
$ws.Range("M134").Value = 1032.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3735.5833
$ws.Range("I31").Value = 1368.875
$ws.Range("K31").Value = 1368.875
$ws.Range("M31").Value = -1073.875
$ws.Range("H34").Value = 3735.5833
$ws.Range("I34").Value = 1368.875
$ws.Range("K34").Value = 1368.875
$ws.Range("M34").Value = -1166.875
$ws.Range("H99").Value = 837039.8
$ws.Range("I99").Value = 1254053
$ws.Range("J99").Value = 3013.5
$ws.Range("K99").Value = 1254053
$ws.Range("L99").Value = 3013.5
$ws.Range("M99").Value = -1252555
$ws.Range("N99").Value = -6009.5
$ws.Range("H126").Value = 837039.8
$ws.Range("I126").Value = 1254053
$ws.Range("J126").Value = 3013.5
$ws.Range("K126").Value = 3762159
$ws.Range("L126").Value = 9040.5
$ws.Range("M126").Value = -3759689
$ws.Range("N126").Value = -13980.5
$ws.Range("H132").Value = 993.8
$ws.Range("I132").Value = 993.8
$ws.Range("K132").Value = 2981.4
$ws.Range("M132").Value = -451.3999999999996
$ws.Range("H134").Value = 3123.8096
$ws.Range("I134").Value = 1431.875
$ws.Range("K134").Value = 4295.625
$ws.Range("M134").Value = -1760.625
$ws.Range("H135").Value = 75166.336
$ws.Range("J135").Value = 75166.336
$ws.Range("L135").Value = 75166.336
$ws.Range("N135").Value = -85306.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 2633.3333
$ws.Range("I76").Value = 2750
$ws.Range("J76").Value = 2400
$ws.Range("K76").Value = 8250
$ws.Range("L76").Value = 7200
$ws.Range("M76").Value = -7867
$ws.Range("N76").Value = -7966
$ws.Range("H79").Value = 2633.3333
$ws.Range("I79").Value = 2750
$ws.Range("J79").Value = 2400
$ws.Range("K79").Value = 8250
$ws.Range("L79").Value = 7200
$ws.Range("M79").Value = -6924
$ws.Range("N79").Value = -9852
$ws.Range("H131").Value = 1753.6666
$ws.Range("I131").Value = 881.125
$ws.Range("J131").Value = 3498.75
$ws.Range("K131").Value = 2643.375
$ws.Range("L131").Value = 10496.25
$ws.Range("M131").Value = 2396.625
$ws.Range("N131").Value = -20576.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1217.9
$ws.Range("I132").Value = 1242.1111
$ws.Range("K132").Value = 3726.3333
$ws.Range("M132").Value = -1196.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5055.353
$ws.Range("I68").Value = 3990.3333
$ws.Range("J68").Value = 5636.273
$ws.Range("K68").Value = 3990.3333
$ws.Range("L68").Value = 5636.273
$ws.Range("M68").Value = -3241.3333
$ws.Range("N68").Value = -7134.273
$ws.Range("H71").Value = 5055.353
$ws.Range("I71").Value = 3990.3333
$ws.Range("J71").Value = 5636.273
$ws.Range("K71").Value = 19951.6665
$ws.Range("L71").Value = 28181.365
$ws.Range("M71").Value = -16207.6665
$ws.Range("N71").Value = -35669.36500000001
$ws.Range("H111").Value = 59500
$ws.Range("J111").Value = 59500
$ws.Range("L111").Value = 59500
$ws.Range("N111").Value = -67680
$ws.Range("H132").Value = 1670.1666
$ws.Range("I132").Value = 1670.1666
$ws.Range("K132").Value = 5010.4998
$ws.Range("M132").Value = -2480.4998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5183.1665
$ws.Range("I62").Value = 5019.8
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 5019.8
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -4395.8
$ws.Range("N62").Value = -7248
$ws.Range("H65").Value = 5183.1665
$ws.Range("I65").Value = 5019.8
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 25099
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -21979
$ws.Range("N65").Value = -36240
$ws.Range("H97").Value = 62000
$ws.Range("J97").Value = 62000
$ws.Range("L97").Value = 62000
$ws.Range("N97").Value = -63982
$ws.Range("H107").Value = 536.1539
$ws.Range("I107").Value = 406.18182
$ws.Range("K107").Value = 1218.54546
$ws.Range("M107").Value = 701.45454
$ws.Range("H132").Value = 1927.6666
$ws.Range("I132").Value = 906.2857
$ws.Range("J132").Value = 5502.5
$ws.Range("K132").Value = 2718.8571
$ws.Range("L132").Value = 16507.5
$ws.Range("M132").Value = -188.8571000000002
$ws.Range("N132").Value = -21567.5
$ws.Range("H136").Value = 2007.5714
$ws.Range("I136").Value = 2007.5714
$ws.Range("K136").Value = 6022.7142
$ws.Range("M136").Value = -3472.7142

Write-Host "Applied 175 cell changes across 8 sheets"
